# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the worker "DUBIS ANTONIA MIELES GALVAN" (row 20) entirely ---
# Deleting the row shifts everything below it up by one (rows 21-31 -> 20-30),
# which also relocates the footer signature rows (old 30/31 -> new 29/30).
$ws.Rows.Item(20).Delete()

# --- Update summary figures at the top of the statement ---
$ws.Range("E11").Value = 245565   # VALOR MORA total
$ws.Range("C13").Value = 5        # Cant. Trabajadores
$ws.Range("F13").Value = 8        # Cant. Periodos

# --- Refresh the worker/period detail table (rows 16-24) with the updated data ---
$ws.Range("B16:G24").ClearContents()

$data = @(
    @(16, "CC", "1047438253", "JORGE ELIECER DIAZ DIAZ",          "1707", 29509, 781242),
    @(17, "CC", "1128059466", "JHAIR MIGUEL PAUTT PEREZ",         "1902", 3312,  877803),
    @(18, "CC", "1143367303", "ALEXANDRA LORDUY GAITAN",          "2001", 12000, 900000),
    @(19, "CC", "1143367303", "ALEXANDRA LORDUY GAITAN",          "2002", 36000, 900000),
    @(20, "CC", "1143367303", "ALEXANDRA LORDUY GAITAN",          "2003", 36000, 900000),
    @(21, "CC", "1042608042", "WENDY LORENA MENDIETA LOPERA",     "2006", 35112, 877803),
    @(22, "CC", "1042608042", "WENDY LORENA MENDIETA LOPERA",     "2007", 35112, 877803),
    @(23, "CC", "1042608042", "WENDY LORENA MENDIETA LOPERA",     "2008", 35112, 877803),
    @(24, "CC", "20269403",   "WILLIAM GUILLERMO FORTICH RINCON", "2008", 23408, 877803)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
